$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69 - this shifts existing rows 69-75 down to 70-76,
# carrying their values/styles with them (matches the diff's net effect of a
# weekly data point being prepended to this market/product series).
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new weekly record.
$ws.Range("A69").Value = 4
$ws.Range("B69").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C69").Value = "Los Lagos"
$ws.Range("D69").Value = 45013
$ws.Range("E69").Value = 10
$ws.Range("F69").Value = 100112030
$ws.Range("G69").Value = "Poroto granado"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 60
$ws.Range("K69").Value = 42000
$ws.Range("L69").Value = 42000
$ws.Range("M69").Value = 42000
$ws.Range("N69").Value = "$/saco 25 kilos"
$ws.Range("O69").Value = "Región Metropolitana"
$ws.Range("P69").Value = 1680
$ws.Range("Q69").Value = 25
$ws.Range("R69").Value = "Hortaliza"

# Make sure the date cell keeps the workbook's date number format,
# same as every other cell in column D.
$ws.Range("D69").NumberFormat = "YYYY-MM-DD HH:MM:SS"
